$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (target raw width 14.42578125 "characters" units;
# the COM ColumnWidth setter here quantizes to 1/6 increments, so
# 13.666667 is the input that lands closest/exactly on the intended width)
$ws.Columns.Item(1).ColumnWidth = 13.666667
$ws.Columns.Item(2).ColumnWidth = 13.666667

# Update cell values
$ws.Range("A1").Value = -0.041771414354548571
$ws.Range("B1").Value = 0.041771414331420578

$ws.Range("A2").Value = 0.048472851300310905
$ws.Range("B2").Value = -0.048472851367908749

$ws.Range("A3").Value = -0.014955321325077672
$ws.Range("B3").Value = 0.01495532119247131
